$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4676
$ws.Range("D2").Value = 80.55

$ws.Range("C3").Value = 1066
$ws.Range("D3").Value = 18.36

$ws.Range("C4").Value = 35
$ws.Range("D4").Value = 0.6

$ws.Range("C5").Value = 28
$ws.Range("D5").Value = 0.48
